$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 139220
$ws.Range("E2").Value = 6446
$ws.Range("F2").Value = 6446
$ws.Range("G2").Value = 7114
$ws.Range("H2").Value = 5362
$ws.Range("I2").Value = 5362
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 61962
$ws.Range("L2").Value = 34251
$ws.Range("M2").Value = 27711
$ws.Range("N2").Value = 27711
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 188
$ws.Range("Q2").Value = 3696
$ws.Range("R2").Value = -4812
$ws.Range("S2").Value = 1506
$ws.Range("T2").Value = 4522
$ws.Range("U2").Value = -826
$ws.Range("V2").Value = 16435
$ws.Range("W2").Value = 4.63
$ws.Range("X2").Value = 3.85
$ws.Range("Y2").Value = 21.01
$ws.Range("Z2").Value = 9.52
$ws.Range("AA2").Value = 123.6
$ws.Range("AB2").Value = 14348.78
$ws.Range("AC2").Value = 14300
$ws.Range("AD2").Value = 20.39
$ws.Range("AE2").Value = 73895
$ws.Range("AF2").Value = 3.94
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 0.6899999999999999
$ws.Range("AI2").Value = 13.99
$ws.Range("AJ2").Value = 37500000

# Row 3
$ws.Range("D3").Value = 146712
$ws.Range("E3").Value = 6980
$ws.Range("F3").Value = 6980
$ws.Range("G3").Value = 5926
$ws.Range("H3").Value = 3768
$ws.Range("I3").Value = 3770
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 74785
$ws.Range("L3").Value = 43817
$ws.Range("M3").Value = 30968
$ws.Range("N3").Value = 30966
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 188
$ws.Range("Q3").Value = 7831
$ws.Range("R3").Value = -8478
$ws.Range("S3").Value = -37
$ws.Range("T3").Value = 3959
$ws.Range("U3").Value = 3872
$ws.Range("V3").Value = 18527
$ws.Range("W3").Value = 4.76
$ws.Range("X3").Value = 2.57
$ws.Range("Y3").Value = 12.85
$ws.Range("Z3").Value = 5.51
$ws.Range("AA3").Value = 141.49
$ws.Range("AB3").Value = 15948.88
$ws.Range("AC3").Value = 10052
$ws.Range("AD3").Value = 19.2
$ws.Range("AE3").Value = 82577
$ws.Range("AF3").Value = 2.34
$ws.Range("AG3").Value = 3000
$ws.Range("AH3").Value = 1.55
$ws.Range("AI3").Value = 29.84
$ws.Range("AJ3").Value = 37500000

# Row 4
$ws.Range("D4").Value = 153406
$ws.Range("E4").Value = 7288
$ws.Range("F4").Value = 7288
$ws.Range("G4").Value = 7039
$ws.Range("H4").Value = 5057
$ws.Range("I4").Value = 5056
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 79675
$ws.Range("L4").Value = 44900
$ws.Range("M4").Value = 34774
$ws.Range("N4").Value = 34772
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 188
$ws.Range("Q4").Value = 6063
$ws.Range("R4").Value = -5175
$ws.Range("S4").Value = -1852
$ws.Range("T4").Value = 3891
$ws.Range("U4").Value = 2172
$ws.Range("V4").Value = 18780
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 3.3
$ws.Range("Y4").Value = 15.38
$ws.Range("Z4").Value = 6.55
$ws.Range("AA4").Value = 129.12
$ws.Range("AB4").Value = 18053.23
$ws.Range("AC4").Value = 13483
$ws.Range("AD4").Value = 11.46
$ws.Range("AE4").Value = 92725
$ws.Range("AF4").Value = 1.67
$ws.Range("AG4").Value = 3000
$ws.Range("AH4").Value = 1.94
$ws.Range("AI4").Value = 22.25
$ws.Range("AJ4").Value = 37500000

# Row 5
$ws.Range("D5").Value = 163583
$ws.Range("E5").Value = 7271
$ws.Range("F5").Value = 7271
$ws.Range("G5").Value = 8893
$ws.Range("H5").Value = 6805
$ws.Range("I5").Value = 6805
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 81857
$ws.Range("L5").Value = 41961
$ws.Range("M5").Value = 39895
$ws.Range("N5").Value = 39893
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 188
$ws.Range("Q5").Value = 6065
$ws.Range("R5").Value = -3449
$ws.Range("S5").Value = -2977
$ws.Range("T5").Value = 3675
$ws.Range("U5").Value = 2391
$ws.Range("V5").Value = 15989
$ws.Range("W5").Value = 4.45
$ws.Range("X5").Value = 4.16
$ws.Range("Y5").Value = 18.23
$ws.Range("Z5").Value = 8.43
$ws.Range("AA5").Value = 105.18
$ws.Range("AB5").Value = 21105.97
$ws.Range("AC5").Value = 18147
$ws.Range("AD5").Value = 7.49
$ws.Range("AE5").Value = 106381
$ws.Range("AF5").Value = 1.28
$ws.Range("AG5").Value = 3000
$ws.Range("AH5").Value = 2.21
$ws.Range("AI5").Value = 16.53
$ws.Range("AJ5").Value = 37500000

# Row 6
$ws.Range("D6").Value = 168656
$ws.Range("E6").Value = 7101
$ws.Range("F6").Value = 7101
$ws.Range("G6").Value = 6263
$ws.Range("H6").Value = 4374
$ws.Range("I6").Value = 4374
$ws.Range("K6").Value = 86878
$ws.Range("L6").Value = 44277
$ws.Range("M6").Value = 42602
$ws.Range("N6").Value = 42602
$ws.Range("P6").Value = 188
$ws.Range("Q6").Value = 4897
$ws.Range("R6").Value = -3895
$ws.Range("S6").Value = 76
$ws.Range("T6").Value = 1812
$ws.Range("U6").Value = 3084
$ws.Range("V6").Value = 18546
$ws.Range("W6").Value = 4.21
$ws.Range("X6").Value = 2.59
$ws.Range("Y6").Value = 10.6
$ws.Range("Z6").Value = 5.18
$ws.Range("AA6").Value = 103.93
$ws.Range("AB6").Value = 22778.49
$ws.Range("AC6").Value = 11663
$ws.Range("AD6").Value = 11.06
$ws.Range("AE6").Value = 113605
$ws.Range("AF6").Value = 1.14
$ws.Range("AG6").Value = 3300
$ws.Range("AH6").Value = 2.56
$ws.Range("AI6").Value = 28.29
$ws.Range("AJ6").Value = 37500000

# Row 7
$ws.Range("D7").Value = 182049
$ws.Range("E7").Value = 8674
$ws.Range("G7").Value = 7133
$ws.Range("H7").Value = 4685
$ws.Range("I7").Value = 4719
$ws.Range("K7").Value = 97239
$ws.Range("L7").Value = 50922
$ws.Range("M7").Value = 46317
$ws.Range("N7").Value = 46231
$ws.Range("P7").Value = 189
$ws.Range("Q7").Value = 6450
$ws.Range("R7").Value = -3444
$ws.Range("S7").Value = -1405
$ws.Range("T7").Value = 3505
$ws.Range("U7").Value = 4548
$ws.Range("W7").Value = 4.76
$ws.Range("X7").Value = 2.57
$ws.Range("Y7").Value = 10.63
$ws.Range("Z7").Value = 5.09
$ws.Range("AA7").Value = 109.94
$ws.Range("AC7").Value = 12585
$ws.Range("AD7").Value = 11.48
$ws.Range("AE7").Value = 123282
$ws.Range("AF7").Value = 1.17
$ws.Range("AG7").Value = 3212
$ws.Range("AH7").Value = 2.22
$ws.Range("AI7").Value = 25.53

# Row 8
$ws.Range("D8").Value = 190325
$ws.Range("E8").Value = 9082
$ws.Range("G8").Value = 8607
$ws.Range("H8").Value = 6100
$ws.Range("I8").Value = 6100
$ws.Range("K8").Value = 103848
$ws.Range("L8").Value = 52470
$ws.Range("M8").Value = 51379
$ws.Range("N8").Value = 51298
$ws.Range("P8").Value = 189
$ws.Range("Q8").Value = 7181
$ws.Range("R8").Value = -4786
$ws.Range("S8").Value = -1649
$ws.Range("T8").Value = 4151
$ws.Range("U8").Value = 3652
$ws.Range("W8").Value = 4.77
$ws.Range("X8").Value = 3.21
$ws.Range("Y8").Value = 12.51
$ws.Range("Z8").Value = 6.07
$ws.Range("AA8").Value = 102.12
$ws.Range("AC8").Value = 16266
$ws.Range("AD8").Value = 9.039999999999999
$ws.Range("AE8").Value = 136794
$ws.Range("AF8").Value = 1.07
$ws.Range("AG8").Value = 3343
$ws.Range("AH8").Value = 2.27
$ws.Range("AI8").Value = 20.55

# Row 9
$ws.Range("D9").Value = 198180
$ws.Range("E9").Value = 9533
$ws.Range("G9").Value = 9134
$ws.Range("H9").Value = 6464
$ws.Range("I9").Value = 6464
$ws.Range("K9").Value = 109624
$ws.Range("L9").Value = 53038
$ws.Range("M9").Value = 56586
$ws.Range("N9").Value = 56447
$ws.Range("P9").Value = 189
$ws.Range("Q9").Value = 7581
$ws.Range("R9").Value = -4927
$ws.Range("S9").Value = -1448
$ws.Range("T9").Value = 4415
$ws.Range("U9").Value = 3826
$ws.Range("W9").Value = 4.81
$ws.Range("X9").Value = 3.26
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 6.06
$ws.Range("AA9").Value = 93.73
$ws.Range("AC9").Value = 17236
$ws.Range("AD9").Value = 8.529999999999999
$ws.Range("AE9").Value = 150527
$ws.Range("AF9").Value = 0.98
$ws.Range("AG9").Value = 3489
$ws.Range("AH9").Value = 2.37
$ws.Range("AI9").Value = 20.24

Write-Output "ok: applied cell updates"
